# Re-generated experiment task-order sheets.
# Each worksheet keeps its position (sheet1.xml..sheet5.xml / rId1..rId5) but
# gets a new name and a freshly "randomised" set of stimulus-file rows, and
# the number of data rows per sheet changes (rows are added/removed as
# needed while re-using the existing header/index-column formatting).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 (rId1): GNG_TO-... -> vSAT_TO-...  (row count unchanged: 4)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "vSAT_TO-16515889440387332"
$ws1.Range("B2").Value = "vSAT_stims-16515889440235972.csv"
$ws1.Range("B3").Value = "SAT_stims-16515889439925416.csv"
$ws1.Range("B4").Value = "SAT_stims-1651588943970868.csv"
$ws1.Range("B5").Value = "vSAT_stims-16515889440084414.csv"

# ---------------------------------------------------------------------
# Sheet 2 (rId2): NB_TO-... -> TOL_TO-...  (row count 9 -> 6, drop rows 8-10)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A8:B10").EntireRow.Delete()
$ws2.Name = "TOL_TO-16515889440865111"
$ws2.Range("B2").Value = "MM_stims-16515889440542455.csv"
$ws2.Range("B3").Value = "ZM_stims-16515889440423322.csv"
$ws2.Range("B4").Value = "MM_stims-16515889440698252.csv"
$ws2.Range("B5").Value = "ZM_stims-16515889440542455.csv"
$ws2.Range("B6").Value = "MM_stims-16515889440855525.csv"
$ws2.Range("B7").Value = "ZM_stims-165158894407185.csv"

# ---------------------------------------------------------------------
# Sheet 3 (rId3): RS_TO-... -> NB_TO-...  (row count 2 -> 9, add rows 4-10)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A3").Copy()
$ws3.Range("A4:A10").PasteSpecial(-4122)
$ws3.Range("A4").Value = 2
$ws3.Range("A5").Value = 3
$ws3.Range("A6").Value = 4
$ws3.Range("A7").Value = 5
$ws3.Range("A8").Value = 6
$ws3.Range("A9").Value = 7
$ws3.Range("A10").Value = 8
$ws3.Name = "NB_TO-16515889467918704"
$ws3.Range("B2").Value = "TB-16515889467702308.csv"
$ws3.Range("B3").Value = "OB-16515889447275171.csv"
$ws3.Range("B4").Value = "TB-1651588946532745.csv"
$ws3.Range("B5").Value = "ZB-match_5-16515889443625016.csv"
$ws3.Range("B6").Value = "ZB-match_1-16515889442041507.csv"
$ws3.Range("B7").Value = "OB-16515889449358168.csv"
$ws3.Range("B8").Value = "ZB-match_2-16515889441090124.csv"
$ws3.Range("B9").Value = "TB-16515889464057643.csv"
$ws3.Range("B10").Value = "OB-16515889445998423.csv"

# ---------------------------------------------------------------------
# Sheet 4 (rId4): TOL_TO-... -> RS_TO-...  (row count 6 -> 2, drop rows 4-7)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A4:B7").EntireRow.Delete()
$ws4.Name = "RS_TO-16515889467938523"
$ws4.Range("B2").Value = "eyes open"
$ws4.Range("B3").Value = "eyes closed"

# ---------------------------------------------------------------------
# Sheet 5 (rId5): vSAT_TO-... -> GNG_TO-...  (row count unchanged: 4)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "GNG_TO-16515889468242111"
$ws5.Range("B2").Value = "go_stims-16515889467968955.csv"
$ws5.Range("B3").Value = "GNG_stims-16515889468084018.csv"
$ws5.Range("B4").Value = "go_stims-16515889468094006.csv"
$ws5.Range("B5").Value = "GNG_stims-16515889468232324.csv"
